# Apply the "live property price update" edit:
#  - After the "Completed" bullet that follows the All Ords PE Ratio /
#    marketindex.com.au data-pull note, add a new sentence describing
#    the live monthly data refresh (with "st" superscripted, matching
#    the existing ordinal-day formatting used elsewhere in the doc).
#  - Move the "_GoBack" bookmark from its old location (right after the
#    "Login details" bullet) to its new location inside the sentence we
#    just inserted (right before the final word "month.").

$d = $word.ActiveDocument

# --- Remove the old _GoBack bookmark (currently sits right after the
#     "Login details" bullet) ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBookmark = $d.Bookmarks.Item("_GoBack")
    $oldBookmark.Delete()
}

# --- Find the correct "Completed" bullet. There are many "Completed"
#     bullets in this document, so anchor off the unique hyperlink text
#     that immediately precedes the one we want, then search for
#     "Completed" only after that point. ---
$anchor = $d.Content
$anchor.Find.Execute("marketindex.com.au/data-downloads") | Out-Null

$target = $d.Range($anchor.End, $d.Content.End)
$target.Find.Execute("Completed") | Out-Null

# Collapse to the point right after "Completed" and insert the new sentence.
$insertPoint = $target.End
$target.Collapse(0)
$newSentence = ". New data (if available) will be pulled on the 1st of each month."
$target.InsertAfter($newSentence)
$insertEnd = $insertPoint + $newSentence.Length

# Make the "st" in "1st" superscript, matching the rest of the document's
# ordinal-day formatting (found relative to the freshly inserted text so
# it can't accidentally match elsewhere).
$newRange = $d.Range($insertPoint, $insertEnd)
$newRange.Find.Execute("st") | Out-Null
$newRange.Font.Superscript = $true

# Insert the new _GoBack bookmark right before "month." (i.e. at the end
# of ". New data (if available) will be pulled on the 1st of each ").
$monthRange = $d.Range($insertPoint, $insertEnd)
$monthRange.Find.Execute("month.") | Out-Null
$bookmarkRange = $d.Range($monthRange.Start, $monthRange.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
